$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 19010.875
$ws.Range("I12").Value = 16374.75
$ws.Range("J12").Value = 21647
$ws.Range("K12").Value = 16374.75
$ws.Range("L12").Value = 21647
$ws.Range("M12").Value = -16204.75
$ws.Range("N12").Value = -21987
$ws.Range("H64").Value = 7832.6924
$ws.Range("I64").Value = 7500
$ws.Range("K64").Value = 7500
$ws.Range("M64").Value = -7252
$ws.Range("H67").Value = 7832.6924
$ws.Range("I67").Value = 7500
$ws.Range("K67").Value = 7500
$ws.Range("M67").Value = -6642
$ws.Range("H74").Value = 7393.696
$ws.Range("J74").Value = 7707.524
$ws.Range("L74").Value = 7707.524
$ws.Range("N74").Value = -9579.524000000001
$ws.Range("H76").Value = 7854
$ws.Range("H77").Value = 7393.696
$ws.Range("J77").Value = 7707.524
$ws.Range("L77").Value = 38537.62
$ws.Range("N77").Value = -47897.62
$ws.Range("H79").Value = 7854
$ws.Range("H98").Value = 984.28125
$ws.Range("I98").Value = 822.4838999999999
$ws.Range("K98").Value = 822.4838999999999
$ws.Range("M98").Value = 675.5161000000001
$ws.Range("H99").Value = 1324.1
$ws.Range("I99").Value = 215
$ws.Range("J99").Value = 3912
$ws.Range("K99").Value = 645
$ws.Range("L99").Value = 11736
$ws.Range("M99").Value = 853
$ws.Range("N99").Value = -14732
$ws.Range("H122").Value = 984.28125
$ws.Range("I122").Value = 822.4838999999999
$ws.Range("K122").Value = 2467.4517
$ws.Range("M122").Value = -17.45169999999962
$ws.Range("H132").Value = 2378.44
$ws.Range("I132").Value = 2378.44
$ws.Range("K132").Value = 7135.32
$ws.Range("M132").Value = -4605.32
$ws.Range("H137").Value = 2563.0852
$ws.Range("I137").Value = 2191.1924
$ws.Range("J137").Value = 3023.524
$ws.Range("K137").Value = 6573.5772
$ws.Range("L137").Value = 9070.572
$ws.Range("M137").Value = -4023.5772
$ws.Range("N137").Value = -14170.572

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 646.3158
$ws.Range("I2").Value = 646.3158
$ws.Range("K2").Value = 646.3158
$ws.Range("M2").Value = -533.3158
$ws.Range("H30").Value = 2596.4
$ws.Range("I30").Value = 995.6667
$ws.Range("J30").Value = 4997.5
$ws.Range("K30").Value = 995.6667
$ws.Range("L30").Value = 4997.5
$ws.Range("M30").Value = -845.6667
$ws.Range("N30").Value = -5297.5
$ws.Range("H32").Value = 4912.615
$ws.Range("I32").Value = 3393.2712
$ws.Range("J32").Value = 9630.579
$ws.Range("K32").Value = 3393.2712
$ws.Range("L32").Value = 9630.579
$ws.Range("M32").Value = -3106.2712
$ws.Range("N32").Value = -10204.579
$ws.Range("H74").Value = 75541.64
$ws.Range("I74").Value = 60551.5
$ws.Range("K74").Value = 60551.5
$ws.Range("M74").Value = -59677.5
$ws.Range("H77").Value = 75541.64
$ws.Range("I77").Value = 60551.5
$ws.Range("K77").Value = 302757.5
$ws.Range("M77").Value = -298389.5
$ws.Range("H116").Value = 646.3158
$ws.Range("I116").Value = 646.3158
$ws.Range("K116").Value = 646.3158
$ws.Range("M116").Value = 1647.6842

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 646.3158
$ws.Range("I3").Value = 646.3158
$ws.Range("K3").Value = 646.3158
$ws.Range("M3").Value = -532.3158
$ws.Range("H86").Value = 14023.235
$ws.Range("J86").Value = 31016.334
$ws.Range("L86").Value = 31016.334
$ws.Range("N86").Value = -33262.334
$ws.Range("H89").Value = 14023.235
$ws.Range("J89").Value = 31016.334
$ws.Range("L89").Value = 155081.67
$ws.Range("N89").Value = -166313.67
$ws.Range("H94").Value = 4768.6045
$ws.Range("I94").Value = 704.6129
$ws.Range("J94").Value = 15267.25
$ws.Range("K94").Value = 704.6129
$ws.Range("L94").Value = 15267.25
$ws.Range("M94").Value = -253.6129
$ws.Range("N94").Value = -16169.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3580.543
$ws.Range("I31").Value = 2411.8096
$ws.Range("K31").Value = 2411.8096
$ws.Range("M31").Value = -2116.8096
$ws.Range("H34").Value = 3580.543
$ws.Range("I34").Value = 2411.8096
$ws.Range("K34").Value = 2411.8096
$ws.Range("M34").Value = -2209.8096
$ws.Range("H122").Value = 3544.9524
$ws.Range("I122").Value = 2870
$ws.Range("J122").Value = 4444.8887
$ws.Range("K122").Value = 8610
$ws.Range("L122").Value = 13334.6661
$ws.Range("M122").Value = -6160
$ws.Range("N122").Value = -18234.6661
$ws.Range("H134").Value = 22291.064
$ws.Range("I134").Value = 29345.031
$ws.Range("K134").Value = 88035.09299999999
$ws.Range("M134").Value = -85500.09299999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 26.142857
$ws.Range("I10").Value = 14.5
$ws.Range("J10").Value = 41.666668
$ws.Range("K10").Value = 43.5
$ws.Range("L10").Value = 125.000004
$ws.Range("M10").Value = 95.5
$ws.Range("N10").Value = -403.000004
$ws.Range("H14").Value = 1211.0667
$ws.Range("I14").Value = 1211.0667
$ws.Range("K14").Value = 3633.2001
$ws.Range("M14").Value = -3460.2001
$ws.Range("H33").Value = 20160.4
$ws.Range("I33").Value = 100
$ws.Range("K33").Value = 600
$ws.Range("M33").Value = -317
$ws.Range("H56").Value = 10006329
$ws.Range("I56").Value = 10006329
$ws.Range("K56").Value = 10006329
$ws.Range("M56").Value = -10005799
$ws.Range("H121").Value = 15873760
$ws.Range("I121").Value = 526.75
$ws.Range("J121").Value = 37038070
$ws.Range("K121").Value = 1580.25
$ws.Range("L121").Value = 111114210
$ws.Range("M121").Value = -270.25
$ws.Range("N121").Value = -111116830
$ws.Range("H141").Value = 2670.3333
$ws.Range("I141").Value = 2670.3333
$ws.Range("K141").Value = 8010.999899999999
$ws.Range("M141").Value = -2830.999899999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").ClearContents()
$ws.Range("H102").Value = 41251.96
$ws.Range("I102").Value = 2284.3076
$ws.Range("J102").Value = 80219.62
$ws.Range("K102").Value = 2284.3076
$ws.Range("L102").Value = 80219.62
$ws.Range("M102").Value = -662.3076000000001
$ws.Range("N102").Value = -83463.62
$ws.Range("H132").Value = 4378.44
$ws.Range("I132").Value = 3953.35
$ws.Range("J132").Value = 6078.8
$ws.Range("K132").Value = 11860.05
$ws.Range("L132").Value = 18236.4
$ws.Range("M132").Value = -9330.049999999999
$ws.Range("N132").Value = -23296.4

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 9093.375
$ws.Range("I40").Value = 5749.75
$ws.Range("J40").Value = 12437
$ws.Range("K40").Value = 5749.75
$ws.Range("L40").Value = 12437
$ws.Range("M40").Value = -5613.75
$ws.Range("N40").Value = -12709

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 7999
$ws.Range("I8").Value = 7999
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 7999
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -7859
$ws.Range("N8").ClearContents()
$ws.Range("H62").Value = 5841.9507
$ws.Range("I62").Value = 3322.0588
$ws.Range("J62").Value = 9015.147999999999
$ws.Range("K62").Value = 3322.0588
$ws.Range("L62").Value = 9015.147999999999
$ws.Range("M62").Value = -2698.0588
$ws.Range("N62").Value = -10263.148
$ws.Range("H65").Value = 5841.9507
$ws.Range("I65").Value = 3322.0588
$ws.Range("J65").Value = 9015.147999999999
$ws.Range("K65").Value = 16610.294
$ws.Range("L65").Value = 45075.74
$ws.Range("M65").Value = -13490.294
$ws.Range("N65").Value = -51315.74
$ws.Range("H136").Value = 2235.8696
$ws.Range("I136").Value = 1338.4736
$ws.Range("K136").Value = 4015.4208
$ws.Range("M136").Value = -1465.4208
